$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove workbook protection element
$wb.Unprotect()

$sheetOld = $wb.Worksheets.Item("Sheet")   # sheetId=1, will become "Dados"
$dadosOld = $wb.Worksheets.Item("Dados")   # sheetId=2, data source

# Move the data over to the first sheet (so the surviving sheet keeps sheetId=1)
$dadosOld.UsedRange.Copy($sheetOld.Range("A1"))

# Drop the original "Dados" sheet
$dadosOld.Delete()

# Rename the remaining sheet to "Dados"
$ws = $sheetOld
$ws.Name = "Dados"

# Update cell contents: row 2 login gets shortened, password becomes encrypted token
$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = "gAAAAABoAAhYMTKwlH2j6sPM_MkjzZEVFW3jJTt5B_ajJUjBBD-0PPmpTJWl8N1vyaLBaclIDxLWCRPdnaFim6nCkJygKCBwSQ=="

# Remove rows 3-5 (used to hold sample users a/b/c)
$ws.Rows.Item(3).Resize(3).Delete()

# Widen column B so the token is readable
$ws.Columns.Item(2).ColumnWidth = 109.16666666666667

# Select M6 on the active (tab-selected) sheet, matching the saved view state
$ws.Range("M6").Select()

Write-Host "done"
